$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("H2").Value = 3
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5
$ws.Range("S2").Value = 1.57
$ws.Range("T2").Value = 2.25
$ws.Range("U2").Value = 2.05
$ws.Range("V2").Value = 1.7
$ws.Range("AA2").Value = 23
$ws.Range("AC2").Value = 6.5
$ws.Range("AF2").Value = 67
$ws.Range("AM2").Value = 501
$ws.Range("AO2").Value = 15
$ws.Range("AP2").Value = 29
$ws.Range("AR2").Value = 81
$ws.Range("AT2").Value = 2.25
$ws.Range("AU2").Value = 9
$ws.Range("AV2").Value = 67
$ws.Range("BB2").Value = 301

# Row 4 updates
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("S4").Value = 1.75
$ws.Range("T4").Value = 2.05

# Row 8 updates
$ws.Range("G8").Value = 1.42
$ws.Range("J8").Value = 1.95
$ws.Range("K8").Value = 2.38
$ws.Range("L8").Value = 7
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("W8").Value = 7
$ws.Range("Y8").Value = 9
$ws.Range("AA8").Value = 12
$ws.Range("AD8").Value = 8.5
$ws.Range("AE8").Value = 21
$ws.Range("AX8").Value = 41
$ws.Range("BB8").Value = 351

$wb.Save()
